$wb = $excel.ActiveWorkbook

# Rename "wt" -> "wt_log2_expression" and "dcin5" -> "dcin5_log2_expression"
$wsWt = $wb.Worksheets.Item("wt")
$wsWt.Name = "wt_log2_expression"

$wsDcin5 = $wb.Worksheets.Item("dcin5")
$wsDcin5.Name = "dcin5_log2_expression"

# Activate the wt_log2_expression sheet (moves tabSelected there and updates activeTab)
$wsWt.Activate()

# Update the selection on that sheet to B45
$wsWt.Range("B45").Select()
